# Edit script: apply the documented changes to Memoria.docx
# Strategy: use Range.InsertXML with precise WordprocessingML fragments to
# replace / split / merge paragraphs, matching the target diff exactly.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-XmlPackage {
    param([string]$BodyXml)
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $wNs + '><w:body>' + $BodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# NOTE: this runtime's PowerShell does not bind named (-Param) arguments
# reliably, so New-XmlPackage is always invoked positionally below.

# -------------------------------------------------------------------------
# EDIT 4 (process bottom-most first so paragraph indices above stay valid):
# After the "...RaspiMJPEG" hyperlink paragraph, the first of the three
# trailing empty paragraphs is replaced with the new "Arduino" section.
# -------------------------------------------------------------------------
$arduinoBody = @'
<w:p><w:pPr><w:rPr><w:b/><w:u w:val="single"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:b/><w:u w:val="single"/></w:rPr></w:pPr></w:p><w:p><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:lastRenderedPageBreak/><w:t>Arduino</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r></w:p><w:p><w:r><w:t>1) Se monta el chasis del robot</w:t></w:r></w:p><w:p><w:r><w:t>2) Se crea el programa en arduino para controlar el robot. Seg&#250;n el car&#225;cter que le entre por el puerto serial ejecutar&#225; una de las acciones (avanzar, girar, retroceder o parar)</w:t></w:r><w:r><w:t xml:space="preserve">. La comunicaci&#243;n con los motores se hace mediante una </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>shield</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> de </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>motores</w:t></w:r><w:r><w:t>(</w:t></w:r><w:bookmarkStart w:id="100" w:name="_GoBack"/><w:proofErr w:type="gramEnd"/><w:r><w:t>L9110S</w:t></w:r><w:bookmarkEnd w:id="100"/><w:r><w:t>-H</w:t></w:r><w:r><w:t>)</w:t></w:r><w:r><w:t xml:space="preserve"> conectada a pilas.</w:t></w:r></w:p><w:p><w:r><w:t>3)</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Conexiones en el robot y arduino: </w:t></w:r></w:p><w:p><w:r><w:t>4) Comprobaci&#243;n de distancia: El arduino ejecuta el sensor ultras&#243;nico para saber la distancia cuando el cliente se lo pide al servidor.</w:t></w:r></w:p>
'@

$pArduinoTarget = $d.Paragraphs(29)
$pArduinoTarget.Range.InsertXML((New-XmlPackage $arduinoBody))

# -------------------------------------------------------------------------
# EDIT 3: Of the three consecutive empty bold/underline paragraphs (17,18,19):
#   - remove paragraph 18 entirely (merge its mark away)
#   - paragraph 17 gets the "8) ..." text
#   - paragraph 19 (now 18) stays empty/unchanged
# -------------------------------------------------------------------------
$p17 = $d.Paragraphs(17)
$p18 = $d.Paragraphs(18)
$mergeRange = $d.Range($p17.Range.Start, $p18.Range.Start)
$mergeRange.Delete()

$p8Body = '<w:p w:rsidR="00A24513" w:rsidRDefault="00A24513"><w:pPr><w:rPr><w:b/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:t>8) Se llama a la funci&#243;n de comprobar distancia 60 veces por segundo (cada vuelta de la animaci&#243;n) Pero solo se realizar&#225; la funci&#243;n cada dos segundos gracias a un contador que comprueba que ha sido llamada 120 veces. Tambi&#233;n se puede forzar la llamada. Cuando se llama se actualiza la distancia y la posibilidad de andar del veh&#237;culo. Si la llamada sale mal se fuerza otra llamada.</w:t></w:r></w:p>'
$p17 = $d.Paragraphs(17)
$p17.Range.InsertXML((New-XmlPackage $p8Body))

# -------------------------------------------------------------------------
# EDIT 2: The empty paragraph 16 (right after "...de CSS.") is replaced by
# two new paragraphs containing points 6) and 7).
# -------------------------------------------------------------------------
$p67Body = '<w:p><w:r><w:t>6) Se completa la API con el robot (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>controlRobot.php</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> y Robot.js) con las instrucciones de acelerar, parar, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>girarIz</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>girarDe</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> y retroceder. Estas son llamadas mediante Ajax seg&#250;n el estado del personaje.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">7) Se a&#241;ade la sombra del personaje. Para ello se utiliza un plano con </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ShadowMaterial</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">.  Adem&#225;s hay que: activar sombras en </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>render</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, poner plano como recibidor de sombras y poner personaje como creador de sombras, adem&#225;s de crear una luz y posicionarla.</w:t></w:r></w:p>'
$p16 = $d.Paragraphs(16)
$p16.Range.InsertXML((New-XmlPackage $p67Body))

# -------------------------------------------------------------------------
# EDIT 1: Remove the bookmarkStart/bookmarkEnd for "_GoBack" inside
# paragraph 15 (the "5) Cuando estemos con un m?vil..." paragraph).
# -------------------------------------------------------------------------
$p15Body = '<w:p w:rsidR="004F549B" w:rsidRDefault="004F549B" w:rsidP="007B5CD4"><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">5) Cuando estemos con un m&#243;vil aparecer&#225;n controles </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>alternativos</w:t></w:r><w:r w:rsidR="003115B9"><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="003115B9"><w:t>flechas  vsg)</w:t></w:r><w:r><w:t xml:space="preserve">. Para saber si estamos con un m&#243;vil se usar&#225; media </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>query</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> de CSS.</w:t></w:r></w:p>'
$p15 = $d.Paragraphs(15)
$p15.Range.InsertXML((New-XmlPackage $p15Body))
